$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the stored timestamp precision in A11 (tiny float drift from a recalculated task run)
$ws.Range("A11").Value = 45866.50026385417

# Append the new row of scheduled-task telemetry (row 12)
$ws.Range("A12").Value = 45866.54189235611
$ws.Range("B12").Value = 2025
$ws.Range("C12").Value = 31
$ws.Range("D12").Value = 20.89
$ws.Range("E12").Value = 70.23999999999999
$ws.Range("F12").Value = 618.75
$ws.Range("G12").Value = 13.34
$ws.Range("H12").Value = "ESE"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "13:00:19"

# Match the existing date-stamp style used by column A (style id 2 -> numFmt 165)
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat
